# Linyola - Matrimonis (Mormons): afegeix la indexacio de cognoms 1852-1856
# (files 96-143), completant el llibre amb les dades del rotlle SPN 2,02 C / A,5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columnes constants de les noves files: es fixen com a rangs abans que les
# dades de cada fila perque l'ordre de creacio de cadenes compartides coincideixi
# amb el del llibre original (Anys "1852-1856", Serie "A,5", Projecte i Rollo).
$ws.Range("M96:M143").Value = '1852-1856'
$ws.Range("L96:L143").Value = 'A,5'
$ws.Range("J96:J143").Value = 'SPN 2,02 C'
$ws.Range("K96:K143").Value = 47

# Any (A) i Fotograma (B) per cada nova fila
$ws.Range("A96").Value = 1852
$ws.Range("B96").Value = 4
$ws.Range("A97").Value = 1852
$ws.Range("B97").Value = 5
$ws.Range("A98").Value = 1852
$ws.Range("B98").Value = 6
$ws.Range("A99").Value = 1852
$ws.Range("B99").Value = 7
$ws.Range("A100").Value = 1852
$ws.Range("B100").Value = 8
$ws.Range("A101").Value = 1852
$ws.Range("B101").Value = 9
$ws.Range("A102").Value = 1852
$ws.Range("B102").Value = 10
$ws.Range("A103").Value = 1852
$ws.Range("B103").Value = 11
$ws.Range("A104").Value = 1853
$ws.Range("B104").Value = 12
$ws.Range("A105").Value = 1853
$ws.Range("B105").Value = 13
$ws.Range("A106").Value = 1853
$ws.Range("B106").Value = 14
$ws.Range("A107").Value = 1853
$ws.Range("B107").Value = 15
$ws.Range("A108").Value = 1853
$ws.Range("B108").Value = 16
$ws.Range("A109").Value = 1853
$ws.Range("B109").Value = 17
$ws.Range("A110").Value = 1854
$ws.Range("B110").Value = 19
$ws.Range("A111").Value = 1854
$ws.Range("B111").Value = 20
$ws.Range("A112").Value = 1854
$ws.Range("B112").Value = 21
$ws.Range("A113").Value = 1854
$ws.Range("B113").Value = 22
$ws.Range("A114").Value = 1854
$ws.Range("B114").Value = 23
$ws.Range("A115").Value = 1854
$ws.Range("B115").Value = 24
$ws.Range("A116").Value = 1854
$ws.Range("B116").Value = 25
$ws.Range("A117").Value = 1854
$ws.Range("B117").Value = 26
$ws.Range("A118").Value = 1854
$ws.Range("B118").Value = 27
$ws.Range("A119").Value = 1855
$ws.Range("B119").Value = 28
$ws.Range("A120").Value = 1855
$ws.Range("B120").Value = 29
$ws.Range("A121").Value = 1855
$ws.Range("B121").Value = 30
$ws.Range("A122").Value = 1855
$ws.Range("B122").Value = 31
$ws.Range("A123").Value = 1855
$ws.Range("B123").Value = 32
$ws.Range("A124").Value = 1855
$ws.Range("B124").Value = 33
$ws.Range("A125").Value = 1855
$ws.Range("B125").Value = 34
$ws.Range("A126").Value = 1855
$ws.Range("B126").Value = 35
$ws.Range("A127").Value = 1855
$ws.Range("B127").Value = 36
$ws.Range("A128").Value = 1855
$ws.Range("B128").Value = 37
$ws.Range("A129").Value = 1855
$ws.Range("B129").Value = 38
$ws.Range("A130").Value = 1855
$ws.Range("B130").Value = 39
$ws.Range("A131").Value = 1856
$ws.Range("B131").Value = 40
$ws.Range("A132").Value = 1856
$ws.Range("B132").Value = 41
$ws.Range("A133").Value = 1856
$ws.Range("B133").Value = 41
$ws.Range("A134").Value = 1856
$ws.Range("B134").Value = 42
$ws.Range("A135").Value = 1856
$ws.Range("B135").Value = 42
$ws.Range("A136").Value = 1856
$ws.Range("B136").Value = 43
$ws.Range("A137").Value = 1856
$ws.Range("B137").Value = 43
$ws.Range("A138").Value = 1856
$ws.Range("B138").Value = 44
$ws.Range("A139").Value = 1856
$ws.Range("B139").Value = 44
$ws.Range("A140").Value = 1856
$ws.Range("B140").Value = 45
$ws.Range("A141").Value = 1856
$ws.Range("B141").Value = 45
$ws.Range("A142").Value = 1856
$ws.Range("B142").Value = 46
$ws.Range("A143").Value = 1856
$ws.Range("B143").Value = 46

# Cognoms Familia (C), en ordre, per cada nova fila
$ws.Range("C96").Value = 'Pasqual Gilavert'
$ws.Range("C97").Value = 'Majoral Fiquera'
$ws.Range("C98").Value = 'Majoral Vergé'
$ws.Range("C99").Value = 'Gormiguera Ginestà'
$ws.Range("C100").Value = 'Formiguera Martí'
$ws.Range("C101").Value = 'Varniol Castelló'
$ws.Range("C102").Value = 'Prenafeta Torres'
$ws.Range("C103").Value = 'Gomà Galceran'
$ws.Range("C104").Value = 'Cava Codina'
$ws.Range("C105").Value = 'Vilaplana Cercós'
$ws.Range("C106").Value = 'Oliva Trepat'
$ws.Range("C107").Value = 'Cercós Martí'
$ws.Range("C108").Value = 'Centena Mata'
$ws.Range("C109").Value = 'Ribes Martí'
$ws.Range("C110").Value = 'Palou Mas'
$ws.Range("C111").Value = 'Palou Mosset'
$ws.Range("C112").Value = 'Bellet Pedrós'
$ws.Range("C113").Value = 'Labaquial Figuera'
$ws.Range("C114").Value = 'Roma Civit'
$ws.Range("C115").Value = 'Binefa Martí'
$ws.Range("C116").Value = 'Formiiguera Bellet'
$ws.Range("C117").Value = 'Fabregat Mosset'
$ws.Range("C118").Value = 'Cascallo Mosset'
$ws.Range("C119").Value = 'Cisteró Nabau'
$ws.Range("C120").Value = 'Riart Rosell'
$ws.Range("C121").Value = 'Roige Pasqual'
$ws.Range("C122").Value = 'Solsona Majoral'
$ws.Range("C123").Value = 'Gili Pujol'
$ws.Range("C124").Value = 'Martí Fradevia'
$ws.Range("C125").Value = 'Moset Mas'
$ws.Range("C126").Value = 'Martí Gené'
$ws.Range("C127").Value = 'Pallerola Solé'
$ws.Range("C128").Value = 'Pedrós Mas'
$ws.Range("C129").Value = 'Cascalló Vallés'
$ws.Range("C130").Value = 'Llovera Sales'
$ws.Range("C131").Value = 'Niubó Mas'
$ws.Range("C132").Value = 'Fabregat Arderiu'
$ws.Range("C133").Value = 'Colell Farré'
$ws.Range("C134").Value = 'Martí Gene'
$ws.Range("C135").Value = 'Civit Martí'
$ws.Range("C136").Value = 'Bellet Escolà'
$ws.Range("C137").Value = 'Renyé Pujol'
$ws.Range("C138").Value = 'Pasqual Ribes'
$ws.Range("C139").Value = 'Mas Riart'
$ws.Range("C140").Value = 'Planes Agulló'
$ws.Range("C141").Value = 'Pedrós Bonet'
$ws.Range("C142").Value = 'Simó Mata'
$ws.Range("C143").Value = 'Torra Civit'

# Deixa la selecio activa a la primera fila buida seguent, com va quedar el llibre
$ws.Range("A144").Select()
